# Applies the "Finalized Experiments with Participant Generation" update:
# renames each task-order sheet with a freshly regenerated timestamp-based
# name, and refreshes the randomized stimulus-file / condition values in
# each sheet's B column (rows 2+) to match the newly generated run.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (new generation timestamps) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16502912424464376"
$wb.Worksheets.Item(2).Name = "NB_TO-1650291246318968"
$wb.Worksheets.Item(3).Name = "RS_TO-16502912463199668"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912463669744"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502912464299724"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912423964.csv"
$ws1.Range("B3").Value = "GNG_stims-1650291242414398.csv"
$ws1.Range("B4").Value = "go_stims-16502912424163957.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912424453948.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_2-16502912427499654.csv"
$ws2.Range("B3").Value = "TB-16502912445508041.csv"
$ws2.Range("B4").Value = "ZB-match_8-1650291242626955.csv"
$ws2.Range("B5").Value = "OB-16502912432999587.csv"
$ws2.Range("B6").Value = "TB-1650291246301984.csv"
$ws2.Range("B7").Value = "OB-16502912433619578.csv"
$ws2.Range("B8").Value = "OB-16502912432349775.csv"
$ws2.Range("B9").Value = "TB-16502912434859638.csv"
$ws2.Range("B10").Value = "ZB-match_2-16502912426659644.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912463339674.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912463219664.csv"
$ws4.Range("B4").Value = "MM_stims-16502912463500023.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912463349664.csv"
$ws4.Range("B6").Value = "MM_stims-165029124636601.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912463509753.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16502912464140072.csv"
$ws5.Range("B3").Value = "SAT_stims-16502912463699718.csv"
$ws5.Range("B4").Value = "vSAT_stims-1650291246398002.csv"
$ws5.Range("B5").Value = "SAT_stims-1650291246382971.csv"
